$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G2").Value = "2016-10-20 09:21:07"

$wsZhCn.Range("H2").Value = "2016-10-20 09:20:55"
$wsZhCn.Range("K2").Value = "2016-10-20 09:21:36"

$wsDeDe.Range("H2").Value = "2016-10-20 09:21:07"
$wsDeDe.Range("K2").Value = "2016-10-20 09:21:54"
